# Auto-generated edit script: updates cryptocurrency price/volume data
# per commit "Updated cryptos list on Thu Dec 28 13:16:04 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.875.43'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '2.383.82'
$ws.Range("E3").Value = '  +4.29%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '''332.55'
$ws.Range("E5").Value = '  +7.39%  '
$ws.Range("D6").Value = '''102.07'
$ws.Range("E6").Value = '  -9.26%  '
$ws.Range("E7").Value = '  +1.48%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +3.22%  '
$ws.Range("D10").Value = '''41.00'
$ws.Range("E10").Value = '  -7.46%  '
$ws.Range("D11").Value = '''0.0932'
$ws.Range("E11").Value = '  +0.49%  '
$ws.Range("D12").Value = '''8.56'
$ws.Range("E12").Value = '  -3.00%  '
$ws.Range("E13").Value = '  -4.42%  '
$ws.Range("D14").Value = '''16.87'
$ws.Range("E14").Value = '  +8.80%  '
$ws.Range("D15").Value = '''0.107'
$ws.Range("E15").Value = '  +1.93%  '
$ws.Range("D16").Value = '2.743.82'
$ws.Range("E16").Value = '  +4.36%  '
$ws.Range("D17").Value = '2.388.13'
$ws.Range("E17").Value = '  +4.55%  '
$ws.Range("D18").Value = '42.929.23'
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").Value = '''7.56'
$ws.Range("E19").Value = '  +5.23%  '
$ws.Range("D20").Value = '''0.0000107'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("E21").Value = '  +8.42%  '
$ws.Range("D22").Value = '''76.39'
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").Value = '''271.59'
$ws.Range("E23").Value = '  +5.72%  '
$ws.Range("E24").Value = '  -3.05%  '
$ws.Range("D25").Value = '''9.86'
$ws.Range("E25").Value = '  +10.10%  '
$ws.Range("D26").Value = '''11.77'
$ws.Range("E26").Value = '  +0.25%  '
$ws.Range("D27").Value = '''1.00'
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '''24.14'
$ws.Range("E28").Value = '  +8.15%  '
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("D30").Value = '''173.86'
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("D31").Value = '''3.11'
$ws.Range("E31").Value = '  -2.33%  '
$ws.Range("D32").Value = '''36.39'
$ws.Range("E32").Value = '  -6.13%  '
$ws.Range("D33").Value = '''0.0921'
$ws.Range("E33").Value = '  +2.41%  '
$ws.Range("D34").Value = '''6.03'
$ws.Range("E34").Value = '  +5.62%  '
$ws.Range("E35").Value = '  +3.67%  '
$ws.Range("E36").Value = '  -5.19%  '
$ws.Range("E37").Value = '  -3.89%  '
$ws.Range("E38").Value = '  -3.62%  '
$ws.Range("D39").Value = '''0.107'
$ws.Range("E39").Value = '  +3.27%  '
$ws.Range("D40").Value = '''2.84'
$ws.Range("E40").Value = '  +12.34%  '
$ws.Range("D41").Value = '''1.54'
$ws.Range("E41").Value = '  +10.64%  '
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("D43").Value = '''69.91'
$ws.Range("E43").Value = '  -3.39%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '''92.86'
$ws.Range("E45").Value = '  +45.24%  '
$ws.Range("D46").Value = '''117.76'
$ws.Range("E46").Value = '  +8.94%  '
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").Value = '''5.51'
$ws.Range("E48").Value = '  -3.37%  '
$ws.Range("D49").Value = '''9.12'
$ws.Range("E49").Value = '  +2.60%  '
$ws.Range("D50").Value = '1.626.67'
$ws.Range("E50").Value = '  +9.74%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = '''0.492'
$ws.Range("E51").Value = '  +13.06%  '
